$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve original text formatting for Price column (avoid Excel
# auto-converting numeric-looking strings like "1.00" or "0.0405" into numbers).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '67.682.29'
$ws.Range("E2").Value = '  +0.16%  '
$ws.Range("D3").Value = '3.331.77'
$ws.Range("E3").Value = '  +0.63%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").Value = '580.38'
$ws.Range("E5").Value = '  -0.99%  '
$ws.Range("D6").Value = '175.68'
$ws.Range("E6").Value = '  -3.35%  '
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '0.587'
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("D9").Value = '3.327.57'
$ws.Range("E9").Value = '  +0.54%  '
$ws.Range("E10").Value = '  +0.88%  '
$ws.Range("D11").Value = '0.576'
$ws.Range("E11").Value = '  -0.04%  '
$ws.Range("D12").Value = '45.36'
$ws.Range("E12").Value = '  -1.69%  '
$ws.Range("E13").Value = '  -1.60%  '
$ws.Range("D14").Value = '669.90'
$ws.Range("E14").Value = '  +5.59%  '
$ws.Range("D15").Value = '3.878.10'
$ws.Range("E15").Value = '  +0.53%  '
$ws.Range("D16").Value = '8.40'
$ws.Range("E16").Value = '  -0.07%  '
$ws.Range("D17").Value = '67.694.10'
$ws.Range("E17").Value = '  -0.13%  '
$ws.Range("E18").Value = '  -0.63%  '
$ws.Range("D19").Value = '3.335.97'
$ws.Range("E19").Value = '  +0.38%  '
$ws.Range("D20").Value = '17.39'
$ws.Range("E20").Value = '  -1.40%  '
$ws.Range("D21").Value = '10.95'
$ws.Range("E21").Value = '  +0.75%  '
$ws.Range("D22").Value = '0.888'
$ws.Range("E22").Value = '  -1.08%  '
$ws.Range("D23").Value = '5.42'
$ws.Range("E23").Value = '  +8.46%  '
$ws.Range("D24").Value = '17.09'
$ws.Range("E24").Value = '  -2.92%  '
$ws.Range("D25").Value = '98.81'
$ws.Range("E25").Value = '  +1.97%  '
$ws.Range("E26").Value = '  -3.65%  '
$ws.Range("E27").Value = '  -3.49%  '
$ws.Range("D28").Value = '9.26'
$ws.Range("E28").Value = '  -2.99%  '
$ws.Range("D29").Value = '33.58'
$ws.Range("E29").Value = '  +3.13%  '
$ws.Range("D30").Value = '8.41'
$ws.Range("E30").Value = '  -1.30%  '
$ws.Range("D31").Value = '7.33'
$ws.Range("E31").Value = '  +10.65%  '
$ws.Range("D32").Value = '573.26'
$ws.Range("E32").Value = '  -2.94%  '
$ws.Range("D33").Value = '10.95'
$ws.Range("E33").Value = '  +0.44%  '
$ws.Range("E34").Value = '  +1.17%  '
$ws.Range("E35").Value = '  +0.23%  '
$ws.Range("D36").Value = '3.687.08'
$ws.Range("E36").Value = '  -6.07%  '
$ws.Range("D37").Value = '56.65'
$ws.Range("E37").Value = '  +1.75%  '
$ws.Range("D38").Value = '3.33'
$ws.Range("E38").Value = '  -5.35%  '
$ws.Range("D39").Value = '34.35'
$ws.Range("E39").Value = '  +5.69%  '
$ws.Range("E40").Value = '  +0.93%  '
$ws.Range("E41").Value = '  -1.79%  '
$ws.Range("E42").Value = '  -4.27%  '
$ws.Range("E43").Value = '  -1.69%  '
$ws.Range("E44").Value = '  -0.90%  '
$ws.Range("D45").Value = '0.0₃0664'
$ws.Range("E45").Value = '  -2.18%  '
$ws.Range("D46").Value = '0.0405'
$ws.Range("E46").Value = '  -1.65%  '
$ws.Range("D47").Value = '2.60'
$ws.Range("E47").Value = '  +2.28%  '
$ws.Range("E48").Value = '  -0.03%  '
$ws.Range("E49").Value = '  -0.37%  '
$ws.Range("D50").Value = '1.35'
$ws.Range("E50").Value = '  -0.34%  '
$ws.Range("D51").Value = '128.98'
$ws.Range("E51").Value = '  -0.99%  '
